$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Other Night" order list lives in columns I (character) and J (reminder
# text) for rows 5-18. Messiah currently sits first (row 5); it needs to move
# to the end of that block (row 18), with every other entry shifting up one
# row to fill the gap ("messiah night order is after evil").
$firstRow = 5
$lastRow = 18

# Capture the current column I / J values for the block before mutating
# anything, using Value2 so special characters (e.g. "&") round-trip cleanly.
$iVals = @()
$jVals = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $iVals += , $ws.Cells.Item($r, 9).Value2
    $jVals += , $ws.Cells.Item($r, 10).Value2
}

# Rotate the block up by one: old row (firstRow) becomes the new lastRow.
$newIVals = @()
$newJVals = @()
for ($k = 1; $k -lt $iVals.Count; $k++) {
    $newIVals += , $iVals[$k]
    $newJVals += , $jVals[$k]
}
$newIVals += , $iVals[0]
$newJVals += , $jVals[0]

for ($k = 0; $k -lt $newIVals.Count; $k++) {
    $r = $firstRow + $k
    $ws.Cells.Item($r, 9).Value = $newIVals[$k]
    $ws.Cells.Item($r, 10).Value = $newJVals[$k]
}

# Update the visible selection to match the saved view state.
$ws.Range("H18").Select()
